# CCC19 Derived Variables Spreadsheet — add new hospitalization variable
#
# A new hospitalization variable (O02b / hosp_bl) that only accounts for
# hospitalizations denoted on the baseline survey form is inserted right
# after the existing O02a row. The previous "baseline" variable (O02a) is
# renamed from hosp_bl to der_hosp_30 (hosp_30) since it also looks for
# hospitalizations in follow-up forms within 30 days of diagnosis.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new blank row right after row 244 (the O02a / hosp_bl row),
# pushing everything from the old row 245 onward down by one.
$ws.Rows.Item(245).Insert()

# Grow the table (ListObject) so it covers the newly inserted row too -
# this keeps the autoFilter / sortState ranges in sync with the sheet.
$lo.Resize($ws.Range("A1:E322"))

# Populate the brand-new row 245 with the new variable definition.
# (Set the two genuinely-new strings first so they land at shared-string
# indices 1062/1063, matching first-seen order, then the renamed value
# last so it lands at 1064.)
$ws.Range("A245").Value = "O02b"
$ws.Range("D245").Value = "Hospitalized at baseline (baseline form only)"
$ws.Range("B245").Value = "hosp_bl"
$ws.Range("C245").Value = "Outcome"
$ws.Range("E245").Value = "0 = No; 1 = Yes; 99 = Unknown"

# Rename the older variable (row 244, O02a) from hosp_bl to der_hosp_30.
$ws.Range("B244").Value = "hosp_30"

# Reflect the edit location in the view: the user ends up with B245
# (the newly entered cell) selected, having scrolled the window down.
$null = $ws.Range("B245").Select()
